# edit.ps1
# Adds a new "2022-Q3" worksheet (with fund-holding detail data) positioned
# right after "总计", and updates the "总计" summary sheet with a new
# leading row for 2022-Q3 (shifting the existing quarters down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value that should be stored as TEXT (t="inlineStr"/shared
# string), forcing Excel to keep numeric-looking strings (fund codes like
# "000612", percentages like "9.45") as text instead of auto-converting them
# to numbers (which would corrupt leading zeros / change the cell type).
# Non-numeric-looking text is written directly so no extra "quoted text"
# style gets attached to it.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    if ($text -match '^-?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")
$wsQ3 = $wb.Worksheets.Add($null, $sheetTotal)
$wsQ3.Name = "2022-Q3"

# Match the outline settings used by the other sheets in the workbook.
$wsQ3.Outline.SummaryRow = 1
$wsQ3.Outline.SummaryColumn = -4152

# Reference cell carrying the workbook's "label" style (bold, centered,
# thin border) used for header rows and the row-index column.
$styleSrc = $sheetTotal.Range("B1")

# Header row
Set-TextValue $wsQ3.Range("B1") "基金代码"
Set-TextValue $wsQ3.Range("C1") "基金名称"
Set-TextValue $wsQ3.Range("D1") "基金规模"
Set-TextValue $wsQ3.Range("E1") "股票总仓位"
Set-TextValue $wsQ3.Range("F1") "仓位占比"
Set-TextValue $wsQ3.Range("G1") "持有市值(亿元)"
Set-TextValue $wsQ3.Range("H1") "仓位排名"
$styleSrc.Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

$detailRows = @(
    [PSCustomObject]@{ Row=2; A=0; B="163801"; C="中银中国混合（LOF）A"; D="9.45"; E="89.73"; F="7.13"; G="0.6738"; GNum=$false; H=4 },
    [PSCustomObject]@{ Row=3; A=1; B="000612"; C="华宝生态中国混合A"; D="8.33"; E="92.34"; F="6.60"; G="0.5498"; GNum=$false; H=1 },
    [PSCustomObject]@{ Row=4; A=2; B="013247"; C="交银瑞卓三年持有期混合"; D="19.44"; E="67.60"; F="2.19"; G="0.4257"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=5; A=3; B="001118"; C="华宝事件驱动混合"; D="6.22"; E="92.69"; F="6.42"; G="0.3993"; GNum=$false; H=1 },
    [PSCustomObject]@{ Row=6; A=4; B="240004"; C="华宝动力组合混合A"; D="14.13"; E="75.08"; F="2.71"; G="0.3829"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=7; A=5; B="163805"; C="中银动态策略混合A"; D="3.93"; E="92.03"; F="7.29"; G="0.2865"; GNum=$false; H=5 },
    [PSCustomObject]@{ Row=8; A=6; B="000124"; C="华宝服务优选混合"; D="6.11"; E="88.51"; F="4.30"; G="0.2627"; GNum=$false; H=3 },
    [PSCustomObject]@{ Row=9; A=7; B="010114"; C="华宝新兴成长混合"; D="3.18"; E="80.50"; F="7.38"; G="0.2347"; GNum=$false; H=2 },
    [PSCustomObject]@{ Row=10; A=8; B="009411"; C="中银科技创新一年定期开放混合"; D="2.14"; E="94.65"; F="8.19"; G="0.1753"; GNum=$false; H=5 },
    [PSCustomObject]@{ Row=11; A=9; B="163809"; C="中银蓝筹精选灵活配置混合"; D="2.08"; E="79.70"; F="7.13"; G="0.1483"; GNum=$false; H=3 },
    [PSCustomObject]@{ Row=12; A=10; B="010418"; C="财通景气行业混合A"; D="2.72"; E="94.88"; F="4.14"; G="0.1126"; GNum=$false; H=10 },
    [PSCustomObject]@{ Row=13; A=11; B="001088"; C="华宝国策导向混合"; D="2.58"; E="87.58"; F="4.04"; G="0.1042"; GNum=$false; H=6 },
    [PSCustomObject]@{ Row=14; A=12; B="240002"; C="华宝宝康配置混合"; D="4.34"; E="63.59"; F="2.22"; G="0.0963"; GNum=$false; H=8 },
    [PSCustomObject]@{ Row=15; A=13; B="010460"; C="兴业研究精选混合A"; D="2.76"; E="87.68"; F="3.38"; G="0.0933"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=16; A=14; B="009189"; C="华宝成长策略混合"; D="1.62"; E="80.11"; F="5.41"; G="0.0876"; GNum=$false; H=5 },
    [PSCustomObject]@{ Row=17; A=15; B="501015"; C="财通多策略升级混合（LOF）A"; D="2.06"; E="94.80"; F="4.15"; G="0.0855"; GNum=$false; H=10 },
    [PSCustomObject]@{ Row=18; A=16; B="016257"; C="华宝动力组合混合C"; D="2.03"; E="75.08"; F="2.71"; G="0.0550"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=19; A=17; B="005959"; C="财通新视野灵活配置混合C"; D="1.12"; E="94.59"; F="3.99"; G="0.0447"; GNum=$false; H=10 },
    [PSCustomObject]@{ Row=20; A=18; B="001370"; C="中银新趋势灵活配置混合A"; D="2.24"; E="39.15"; F="1.99"; G="0.0446"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=21; A=19; B="002430"; C="中银丰利灵活配置混合A"; D="3.62"; E="20.50"; F="0.99"; G="0.0358"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=22; A=20; B="015271"; C="财通多策略升级混合（LOF）C"; D="0.62"; E="94.80"; F="4.15"; G="0.0257"; GNum=$false; H=10 },
    [PSCustomObject]@{ Row=23; A=21; B="005851"; C="财通新视野灵活配置混合A"; D="0.62"; E="94.59"; F="3.99"; G="0.0247"; GNum=$false; H=10 },
    [PSCustomObject]@{ Row=24; A=22; B="002152"; C="华宝核心优势灵活配置混合A"; D="0.47"; E="91.02"; F="5.22"; G="0.0245"; GNum=$false; H=4 },
    [PSCustomObject]@{ Row=25; A=23; B="009366"; C="浦银安盛科技创新一年定期开放混合A"; D="0.87"; E="89.66"; F="2.30"; G="0.0200"; GNum=$false; H=10 },
    [PSCustomObject]@{ Row=26; A=24; B="016462"; C="华宝生态中国混合C"; D="0.15"; E="92.34"; F="6.60"; G="0.0099"; GNum=$false; H=1 },
    [PSCustomObject]@{ Row=27; A=25; B="011482"; C="中银顺宁回报6个月持有期混合A"; D="0.75"; E="22.09"; F="0.99"; G="0.0074"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=28; A=26; B="519175"; C="浦银安盛经济带崛起灵活配置混合"; D="0.33"; E="34.37"; F="1.51"; G="0.0050"; GNum=$false; H=8 },
    [PSCustomObject]@{ Row=29; A=27; B="007084"; C="天治转型升级混合"; D="0.11"; E="91.86"; F="2.97"; G="0.0033"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=30; A=28; B="004801"; C="浦银安盛安久回报定期开放混合A"; D="0.17"; E="21.44"; F="1.52"; G="0.0026"; GNum=$false; H=4 },
    [PSCustomObject]@{ Row=31; A=29; B="002431"; C="中银丰利灵活配置混合C"; D="0.22"; E="20.50"; F="0.99"; G="0.0022"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=32; A=30; B="008162"; C="浦银安盛经济带崛起灵活配置混合C"; D="0.13"; E="34.37"; F="1.51"; G="0.0020"; GNum=$false; H=8 },
    [PSCustomObject]@{ Row=33; A=31; B="009367"; C="浦银安盛科技创新一年定期开放混合C"; D="0.08"; E="89.66"; F="2.30"; G="0.0018"; GNum=$false; H=10 },
    [PSCustomObject]@{ Row=34; A=32; B="011483"; C="中银顺宁回报6个月持有期混合C"; D="0.08"; E="22.09"; F="0.99"; G="0.0008"; GNum=$false; H=9 },
    [PSCustomObject]@{ Row=35; A=33; B="015365"; C="中银动态策略混合C"; D="0.01"; E="92.03"; F="7.29"; G="0.0007"; GNum=$false; H=5 },
    [PSCustomObject]@{ Row=36; A=34; B="004802"; C="浦银安盛安久回报定期开放混合C"; D="0.00"; E="21.44"; F="1.52"; G="0"; GNum=$true; H=4 },
    [PSCustomObject]@{ Row=37; A=35; B="014537"; C="中银中国混合（LOF）C"; D="0.00"; E="89.73"; F="7.13"; G="0"; GNum=$true; H=4 },
    [PSCustomObject]@{ Row=38; A=36; B="014845"; C="中银新趋势灵活配置混合C"; D="0.00"; E="39.15"; F="1.99"; G="0"; GNum=$true; H=9 },
    [PSCustomObject]@{ Row=39; A=37; B="016234"; C="财通景气行业混合C"; D="0.00"; E="94.88"; F="4.14"; G="0"; GNum=$true; H=10 },
    [PSCustomObject]@{ Row=40; A=38; B="015947"; C="兴业研究精选混合C"; D="0.00"; E="87.68"; F="3.38"; G="0"; GNum=$true; H=9 },
    [PSCustomObject]@{ Row=41; A=39; B="016461"; C="华宝核心优势灵活配置混合C"; D="0.00"; E="91.02"; F="5.22"; G="0"; GNum=$true; H=4 }
)

foreach ($row in $detailRows) {
    $r = $row.Row
    $wsQ3.Range("A$r").Value = $row.A
    Set-TextValue $wsQ3.Range("B$r") $row.B
    Set-TextValue $wsQ3.Range("C$r") $row.C
    Set-TextValue $wsQ3.Range("D$r") $row.D
    Set-TextValue $wsQ3.Range("E$r") $row.E
    Set-TextValue $wsQ3.Range("F$r") $row.F
    if ($row.GNum) {
        $wsQ3.Range("G$r").Value = [double]$row.G
    } else {
        Set-TextValue $wsQ3.Range("G$r") $row.G
    }
    $wsQ3.Range("H$r").Value = $row.H

    $styleSrc.Copy()
    $wsQ3.Range("A$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: add the 2022-Q3 row on top, shifting the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$summaryRows = @(
    [PSCustomObject]@{ Row=2; A=0; B="2022-Q3"; C=40; D=4.43 },
    [PSCustomObject]@{ Row=3; A=1; B="2022-Q2"; C=38; D=8.73 },
    [PSCustomObject]@{ Row=4; A=2; B="2022-Q1"; C=11; D=2.38 },
    [PSCustomObject]@{ Row=5; A=3; B="2021-Q4"; C=21; D=7.27 },
    [PSCustomObject]@{ Row=6; A=4; B="2021-Q3"; C=7; D=2.16 },
    [PSCustomObject]@{ Row=7; A=5; B="2021-Q2"; C=27; D=6.42 },
    [PSCustomObject]@{ Row=8; A=6; B="2021-Q1"; C=29; D=5.51 },
    [PSCustomObject]@{ Row=9; A=7; B="2020-Q4"; C=9; D=11.85 }
)

foreach ($row in $summaryRows) {
    $r = $row.Row
    $sheetTotal.Range("A$r").Value = $row.A
    Set-TextValue $sheetTotal.Range("B$r") $row.B
    $sheetTotal.Range("C$r").Value = $row.C
    $sheetTotal.Range("D$r").Value = $row.D

    $styleSrc.Copy()
    $sheetTotal.Range("A$r").PasteSpecial(-4122)
}

Write-Output "2022-Q3 sheet added and 总计 updated"
